$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.406.81'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '3.098.49'
$ws.Range('E3').Value = '  +2.17%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '387.97'
$ws.Range('E5').Value = '  +1.96%  '
$ws.Range('D6').Value = '103.68'
$ws.Range('E6').Value = '  +0.78%  '
$ws.Range('D7').Value = '0.538'
$ws.Range('E7').Value = '  -1.39%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.588'
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('D10').Value = '37.07'
$ws.Range('E10').Value = '  +0.62%  '
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').Value = '0.0858'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').Value = '3.590.83'
$ws.Range('E13').Value = '  +2.15%  '
$ws.Range('D14').Value = '18.53'
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').Value = '7.78'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').Value = '3.088.38'
$ws.Range('E16').Value = '  +1.73%  '
$ws.Range('D17').Value = '0.998'
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('D18').Value = '10.66'
$ws.Range('E18').Value = '  +1.31%  '
$ws.Range('D19').Value = '51.502.51'
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('D20').Value = '3.25'
$ws.Range('E20').Value = '  +6.05%  '
$ws.Range('D21').Value = '12.52'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').Value = '0.0₃0966'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').Value = '70.21'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').Value = '266.25'
$ws.Range('E24').Value = '  -0.83%  '
$ws.Range('D25').Value = '3.16'
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').Value = '8.02'
$ws.Range('E26').Value = '  -3.26%  '
$ws.Range('D27').Value = '27.40'
$ws.Range('E27').Value = '  +4.28%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '7.17'
$ws.Range('E29').Value = '  -5.99%  '
$ws.Range('E30').Value = '  -4.79%  '
$ws.Range('E31').Value = '  -1.83%  '
$ws.Range('D32').Value = '10.41'
$ws.Range('E32').Value = '  +1.24%  '
$ws.Range('D33').Value = '35.73'
$ws.Range('E33').Value = '  +4.66%  '
$ws.Range('D34').Value = '0.0472'
$ws.Range('E34').Value = '  +5.42%  '
$ws.Range('E35').Value = '  +0.90%  '
$ws.Range('E36').Value = '  -0.99%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').Value = '0.291'
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('D40').Value = '129.50'
$ws.Range('E40').Value = '  +4.44%  '
$ws.Range('E41').Value = '  -0.37%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '0.116'
$ws.Range('E42').Value = '  -0.35%  '
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D43').Value = '16.56'
$ws.Range('E43').Value = '  -3.01%  '
$ws.Range('D44').Value = '3.83'
$ws.Range('E44').Value = '  +1.80%  '
$ws.Range('D45').Value = '2.50'
$ws.Range('E45').Value = '  -2.76%  '
$ws.Range('D46').Value = '22.09'
$ws.Range('E46').Value = '  +1.11%  '
$ws.Range('D47').Value = '2.50'
$ws.Range('E47').Value = '  +4.06%  '
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('D49').Value = '2.076.58'
$ws.Range('E49').Value = '  +1.99%  '
$ws.Range('B50').Value = 'BEAM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D50').Value = '0.0332'
$ws.Range('E50').Value = '  +3.80%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.926'
$ws.Range('E51').Value = '  +18.46%  '
